$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all existing data rows (2-458)
# from 2023-10-06 (45205) to 2023-10-07 (45206)
$ws.Range("C2:C458").Value = 45206

# Row 458 picks up an explicit row height now that a new row follows it
$ws.Rows.Item(458).RowHeight = 15

# Append the new record as row 459
$ws.Cells.Item(459, 1).Value = "A 48081-2023"

$ws.Cells.Item(459, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(459, 2).Value = 45204

$ws.Cells.Item(459, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(459, 3).Value = 45206

$ws.Cells.Item(459, 4).Value = "VÄSTERBOTTENS LÄN"
$ws.Cells.Item(459, 5).Value = "ÅSELE"
$ws.Cells.Item(459, 6).Value = "SCA"
$ws.Cells.Item(459, 7).Value = 0.3
$ws.Cells.Item(459, 8).Value = 0
$ws.Cells.Item(459, 9).Value = 0
$ws.Cells.Item(459, 10).Value = 0
$ws.Cells.Item(459, 11).Value = 0
$ws.Cells.Item(459, 12).Value = 0
$ws.Cells.Item(459, 13).Value = 0
$ws.Cells.Item(459, 14).Value = 0
$ws.Cells.Item(459, 15).Value = 0
$ws.Cells.Item(459, 16).Value = 0
$ws.Cells.Item(459, 17).Value = 0

$ws.Cells.Item(459, 18).WrapText = $true
$ws.Cells.Item(459, 18).Value = ""
